$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bom")

# Set the "y" marker cells first so the shared string "y" is created
# before the new header strings (Mouser/Digikey/Amazon/Ordered?),
# matching the original authoring order.
$orderedRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17)
foreach ($r in $orderedRows) {
    if ($r -eq 13 -or $r -eq 14 -or $r -eq 17) {
        $ws.Cells.Item($r, 9).Value = "y"   # column I
    } else {
        $ws.Cells.Item($r, 7).Value = "y"   # column G
    }
    $ws.Cells.Item($r, 10).Value = "y"      # column J
}

# Header row
$ws.Range("G1").Value = "Mouser"
$ws.Range("H1").Value = "Digikey"
$ws.Range("I1").Value = "Amazon"
$ws.Range("J1").Value = "Ordered?"

# Match bold header style of the existing header row (A1:F1)
$ws.Range("G1:J1").Font.Bold = $true

# Update the active selection as seen in the edited workbook
$ws.Range("P12").Select() | Out-Null
